$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (H1) keeps its own (non-shared) formula, but the number formatting
# inside the CONCATENATE changes: D1 now goes through TEXT(...,"#,##0") for
# thousands-separated display, and E1's scientific notation drops a digit of
# precision (0.00E+00 -> 0.0E+00).
$ws.Range("H1").Formula = '=CONCATENATE("| ",B1," | ",C1," | ",TEXT(D1,"#,##0")," | ",TEXT(E1,"0.0E+00")," | ",TEXT(F1,"0.0E+00")," | ",TEXT(G1,"0.0E+00")," |")'

# Rows 2-24 (H2:H24) get the same updated formula pattern, applied in one
# shot so the engine re-derives the shared-formula block for the range.
$ws.Range("H2:H24").Formula = '=CONCATENATE("| ",B2," | ",C2," | ",TEXT(D2,"#,##0")," | ",TEXT(E2,"0.0E+00")," | ",TEXT(F2,"0.0E+00")," | ",TEXT(G2,"0.0E+00")," |")'
